$d = $word.ActiveDocument

# 1. Update total hours calculated in dataset
$d.Content.Find.Execute(
    "Total time in hours calculated in dataset: 719",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Total time in hours calculated in dataset: 719.92",
    2)

# 2. Update FDD flag 2 total time line (label + value both change)
$d.Content.Find.Execute(
    "Total time for when FDD flag 2 is True: 7 days 18:00:00",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Total time in hours for when FDD flag 2 is True: 18.0",
    2)

# 3. Update FDD flag 3 total time line (label + value both change)
$d.Content.Find.Execute(
    "Total time for when FDD flag 3 is True: 1 days 22:55:00",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Total time in hours for when FDD flag 3 is True: 22.92",
    2)

# 4. Update the "Report generated" timestamp
$d.Content.Find.Execute(
    "Report generated: Sat Dec  3 09:04:54 2022",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Report generated: Mon Dec  5 10:16:17 2022",
    2)

# 5. Remove the empty ListBullet paragraph that sits right before the
#    "Report generated" paragraph.
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    $txt = $p.Range.Text
    if ($txt.Trim() -eq "" -and $p.Style.NameLocal -eq "List Bullet") {
        $next = $p.Next()
        if ($next -ne $null -and $next.Range.Text.StartsWith("Report generated:")) {
            $p.Range.Delete()
            break
        }
    }
}
